$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$flat = @("2", "2.473e-19", "1.602", "0.526", "5", "3", "1.81e-18", "1.995", "0.5928", "5", "4", "0.3462", "1.112", "0.2749", "20", "5", "0.3308", "0.6518", "0.1534", "25", "6", "0.1351", "0.7959000000000001", "0.1706", "25", "7", "0.06888", "1.33", "0.4627", "20", "8", "0.05829", "1.523", "1.399", "15", "9", "0.2477", "1.86", "1.228", "20", "10", "0.4646", "1.368", "0.5923", "25", "11", "0.4228", "1.505", "0.1071", "25", "12", "0.3204", "1.641", "0.7953", "25", "13", "0.002187", "1.546", "0.1693", "10", "14", "0.1009", "2.353", "0.3006", "15", "15", "0.09429999999999999", "2.451", "0.242", "15", "16", "0.00134", "2.062", "0.7192", "10", "17", "0.0002662", "2.281", "0.2886", "10", "18", "0.007861", "2.941", "0.2737", "15", "19", "0.1738", "2.607", "0.2202", "25", "20", "0.01123", "3.101", "0.2173", "15", "21", "0.0126", "3.25", "0.107", "20", "22", "0.008581999999999999", "3.206", "0.1225", "20", "23", "0.03624", "3.094", "0.4257", "25", "24", "0.01635", "3.287", "0.6127", "25", "25", "0.009377", "3.895", "0.6296", "20", "26", "0.01314", "3.782", "0.1966", "20", "27", "0.01183", "4.497", "0.3537", "15", "28", "0.0006742", "4.164", "0.526", "15", "29", "0.02656", "4.59", "0.1488", "20", "30", "0.0008779", "4.464", "0.2842", "15", "31", "0.06807000000000001", "4.294", "0.253", "25", "32", "0.0004758", "3.839", "0.6179", "10", "33", "0.01389", "4.602", "0.2695", "15", "34", "0.005308", "4.643", "0.3045", "15", "35", "0.001326", "4.502", "0.164", "15", "36", "0.0005083", "4.651", "0.2042", "15", "37", "0.07881000000000001", "4.605", "0.2847", "25", "38", "0.06856", "4.639", "0.4637", "25", "39", "0.03904", "4.673", "0.7395", "25", "40", "0.01848", "4.985", "0.2089", "15", "41", "0.006495", "4.989", "0.08316999999999999", "15", "42", "0.003861", "4.937", "0.3786", "15", "43", "0.1434", "4.774", "0.2819", "25", "44", "0.229", "4.86", "1.989", "25", "45", "0.01825", "4.494", "2.31", "10", "46", "0.1221", "5.243", "0.4614", "15", "47", "0.2315", "5.416", "0.6021", "20", "48", "0.406", "5.153", "0.144", "25", "49", "1.754e-18", "6.09", "0.5756", "5", "50", "0.03743", "5.558", "1.615", "15", "51", "1.261e-17", "5.469", "0.9634", "5", "52", "0.02474", "5.327", "1.329", "15", "53", "0.08151", "5.445", "0.3352", "20", "54", "0.6843", "5.093", "0.4235", "25", "55", "0.3566", "5.026", "1.096", "25", "56", "0.07369000000000001", "4.767", "0.1319", "15", "57", "0.02181", "4.609", "0.6167", "15", "58", "0.01284", "3.524", "0.361", "10", "59", "0.243", "4.662", "0.3571", "25", "60", "0.02578", "4.212", "0.07553", "15", "61", "0.194", "4.492", "0.1818", "25", "62", "0.147", "4.398", "0.1987", "25", "63", "0.1728", "4.4", "0.9268999999999999", "25", "64", "0.006758", "3.554", "0.5497", "10", "65", "0.009176", "4.357", "0.2289", "15", "66", "0.002954", "4.328", "1.317", "15", "67", "9.974e-05", "4.331", "0.3696", "10", "68", "0.1508", "4.494", "0.5243", "25", "69", "0.006925", "4.86", "0.1013", "10", "70", "0.02144", "5.413", "0.0873", "15", "71", "0.006183", "5.176", "0.5238", "15", "72", "0.02639", "5.597", "0.3582", "20")
$n = $flat.Count
$i = 0
while ($i -lt $n) {
  $r = [int]$flat[$i]
  $c = [double]$flat[$i+1]
  $d = [double]$flat[$i+2]
  $e = [double]$flat[$i+3]
  $f = [double]$flat[$i+4]
  $ws.Cells.Item($r, 3).Value = $c
  $ws.Cells.Item($r, 4).Value = $d
  $ws.Cells.Item($r, 5).Value = $e
  $ws.Cells.Item($r, 6).Value = $f
  $i = $i + 5
}
